# Replace each arithmetic-problem answer cell's text with its regenerated
# equivalent. All 100 "old" equation strings are unique within the
# document, so a whole-word Find/Replace targeting each one in turn is
# unambiguous and updates exactly the corresponding table cell.
$d = $word.ActiveDocument

$d.Content.Find.Execute("28+48=76", $true, $true, $false, $false, $false, $true, 1, $false, "52-19=33", 2) | Out-Null
$d.Content.Find.Execute("59-18=41", $true, $true, $false, $false, $false, $true, 1, $false, "63-13=50", 2) | Out-Null
$d.Content.Find.Execute("25+3=28", $true, $true, $false, $false, $false, $true, 1, $false, "58+4=62", 2) | Out-Null
$d.Content.Find.Execute("40-1=39", $true, $true, $false, $false, $false, $true, 1, $false, "56-8=48", 2) | Out-Null
$d.Content.Find.Execute("5+4=9", $true, $true, $false, $false, $false, $true, 1, $false, "49-26=23", 2) | Out-Null
$d.Content.Find.Execute("2+39=41", $true, $true, $false, $false, $false, $true, 1, $false, "3+73=76", 2) | Out-Null
$d.Content.Find.Execute("30+23=53", $true, $true, $false, $false, $false, $true, 1, $false, "83-7=76", 2) | Out-Null
$d.Content.Find.Execute("84+6=90", $true, $true, $false, $false, $false, $true, 1, $false, "80-27=53", 2) | Out-Null
$d.Content.Find.Execute("83-9=74", $true, $true, $false, $false, $false, $true, 1, $false, "57-31=26", 2) | Out-Null
$d.Content.Find.Execute("46-32=14", $true, $true, $false, $false, $false, $true, 1, $false, "69-21=48", 2) | Out-Null
$d.Content.Find.Execute("51-26=25", $true, $true, $false, $false, $false, $true, 1, $false, "25-5=20", 2) | Out-Null
$d.Content.Find.Execute("39-12=27", $true, $true, $false, $false, $false, $true, 1, $false, "75-8=67", 2) | Out-Null
$d.Content.Find.Execute("34+16=50", $true, $true, $false, $false, $false, $true, 1, $false, "51-23=28", 2) | Out-Null
$d.Content.Find.Execute("6+0=6", $true, $true, $false, $false, $false, $true, 1, $false, "83-42=41", 2) | Out-Null
$d.Content.Find.Execute("62+17=79", $true, $true, $false, $false, $false, $true, 1, $false, "62+15=77", 2) | Out-Null
$d.Content.Find.Execute("70-64=6", $true, $true, $false, $false, $false, $true, 1, $false, "27+38=65", 2) | Out-Null
$d.Content.Find.Execute("33-7=26", $true, $true, $false, $false, $false, $true, 1, $false, "67-29=38", 2) | Out-Null
$d.Content.Find.Execute("40+36=76", $true, $true, $false, $false, $false, $true, 1, $false, "46+53=99", 2) | Out-Null
$d.Content.Find.Execute("51-33=18", $true, $true, $false, $false, $false, $true, 1, $false, "25+24=49", 2) | Out-Null
$d.Content.Find.Execute("14+79=93", $true, $true, $false, $false, $false, $true, 1, $false, "24-7=17", 2) | Out-Null
$d.Content.Find.Execute("13+12=25", $true, $true, $false, $false, $false, $true, 1, $false, "47+33=80", 2) | Out-Null
$d.Content.Find.Execute("56+34=90", $true, $true, $false, $false, $false, $true, 1, $false, "50-2=48", 2) | Out-Null
$d.Content.Find.Execute("23+7=30", $true, $true, $false, $false, $false, $true, 1, $false, "49-27=22", 2) | Out-Null
$d.Content.Find.Execute("34+48=82", $true, $true, $false, $false, $false, $true, 1, $false, "36+34=70", 2) | Out-Null
$d.Content.Find.Execute("76-51=25", $true, $true, $false, $false, $false, $true, 1, $false, "48-23=25", 2) | Out-Null
$d.Content.Find.Execute("39+10=49", $true, $true, $false, $false, $false, $true, 1, $false, "93-88=5", 2) | Out-Null
$d.Content.Find.Execute("51+26=77", $true, $true, $false, $false, $false, $true, 1, $false, "15+48=63", 2) | Out-Null
$d.Content.Find.Execute("78-36=42", $true, $true, $false, $false, $false, $true, 1, $false, "93-52=41", 2) | Out-Null
$d.Content.Find.Execute("38+61=99", $true, $true, $false, $false, $false, $true, 1, $false, "15+60=75", 2) | Out-Null
$d.Content.Find.Execute("90-19=71", $true, $true, $false, $false, $false, $true, 1, $false, "56-4=52", 2) | Out-Null
$d.Content.Find.Execute("39-35=4", $true, $true, $false, $false, $false, $true, 1, $false, "76-37=39", 2) | Out-Null
$d.Content.Find.Execute("52+15=67", $true, $true, $false, $false, $false, $true, 1, $false, "68-38=30", 2) | Out-Null
$d.Content.Find.Execute("65-13=52", $true, $true, $false, $false, $false, $true, 1, $false, "98-31=67", 2) | Out-Null
$d.Content.Find.Execute("79-56=23", $true, $true, $false, $false, $false, $true, 1, $false, "91-61=30", 2) | Out-Null
$d.Content.Find.Execute("63-49=14", $true, $true, $false, $false, $false, $true, 1, $false, "47+22=69", 2) | Out-Null
$d.Content.Find.Execute("3+88=91", $true, $true, $false, $false, $false, $true, 1, $false, "9+9=18", 2) | Out-Null
$d.Content.Find.Execute("62-3=59", $true, $true, $false, $false, $false, $true, 1, $false, "90-18=72", 2) | Out-Null
$d.Content.Find.Execute("19+71=90", $true, $true, $false, $false, $false, $true, 1, $false, "98-3=95", 2) | Out-Null
$d.Content.Find.Execute("58+5=63", $true, $true, $false, $false, $false, $true, 1, $false, "83-43=40", 2) | Out-Null
$d.Content.Find.Execute("62-42=20", $true, $true, $false, $false, $false, $true, 1, $false, "44+11=55", 2) | Out-Null
$d.Content.Find.Execute("14-9=5", $true, $true, $false, $false, $false, $true, 1, $false, "80-9=71", 2) | Out-Null
$d.Content.Find.Execute("13+83=96", $true, $true, $false, $false, $false, $true, 1, $false, "76-42=34", 2) | Out-Null
$d.Content.Find.Execute("75-35=40", $true, $true, $false, $false, $false, $true, 1, $false, "81+1=82", 2) | Out-Null
$d.Content.Find.Execute("74-11=63", $true, $true, $false, $false, $false, $true, 1, $false, "77-13=64", 2) | Out-Null
$d.Content.Find.Execute("88-75=13", $true, $true, $false, $false, $false, $true, 1, $false, "49-19=30", 2) | Out-Null
$d.Content.Find.Execute("28+6=34", $true, $true, $false, $false, $false, $true, 1, $false, "61-29=32", 2) | Out-Null
$d.Content.Find.Execute("90-88=2", $true, $true, $false, $false, $false, $true, 1, $false, "14+27=41", 2) | Out-Null
$d.Content.Find.Execute("25+57=82", $true, $true, $false, $false, $false, $true, 1, $false, "59+32=91", 2) | Out-Null
$d.Content.Find.Execute("43-3=40", $true, $true, $false, $false, $false, $true, 1, $false, "36-6=30", 2) | Out-Null
$d.Content.Find.Execute("43+16=59", $true, $true, $false, $false, $false, $true, 1, $false, "16-12=4", 2) | Out-Null
$d.Content.Find.Execute("16+14=30", $true, $true, $false, $false, $false, $true, 1, $false, "19+46=65", 2) | Out-Null
$d.Content.Find.Execute("91-33=58", $true, $true, $false, $false, $false, $true, 1, $false, "63-24=39", 2) | Out-Null
$d.Content.Find.Execute("39-22=17", $true, $true, $false, $false, $false, $true, 1, $false, "38+34=72", 2) | Out-Null
$d.Content.Find.Execute("11+49=60", $true, $true, $false, $false, $false, $true, 1, $false, "43-39=4", 2) | Out-Null
$d.Content.Find.Execute("20-15=5", $true, $true, $false, $false, $false, $true, 1, $false, "50+37=87", 2) | Out-Null
$d.Content.Find.Execute("10+26=36", $true, $true, $false, $false, $false, $true, 1, $false, "48+33=81", 2) | Out-Null
$d.Content.Find.Execute("38+43=81", $true, $true, $false, $false, $false, $true, 1, $false, "23+29=52", 2) | Out-Null
$d.Content.Find.Execute("85-30=55", $true, $true, $false, $false, $false, $true, 1, $false, "59-47=12", 2) | Out-Null
$d.Content.Find.Execute("91-68=23", $true, $true, $false, $false, $false, $true, 1, $false, "26-23=3", 2) | Out-Null
$d.Content.Find.Execute("32+0=32", $true, $true, $false, $false, $false, $true, 1, $false, "99-51=48", 2) | Out-Null
$d.Content.Find.Execute("73+1=74", $true, $true, $false, $false, $false, $true, 1, $false, "96-49=47", 2) | Out-Null
$d.Content.Find.Execute("8+81=89", $true, $true, $false, $false, $false, $true, 1, $false, "2+62=64", 2) | Out-Null
$d.Content.Find.Execute("81-12=69", $true, $true, $false, $false, $false, $true, 1, $false, "0+69=69", 2) | Out-Null
$d.Content.Find.Execute("87-14=73", $true, $true, $false, $false, $false, $true, 1, $false, "90-43=47", 2) | Out-Null
$d.Content.Find.Execute("72+11=83", $true, $true, $false, $false, $false, $true, 1, $false, "18+39=57", 2) | Out-Null
$d.Content.Find.Execute("35+9=44", $true, $true, $false, $false, $false, $true, 1, $false, "55-46=9", 2) | Out-Null
$d.Content.Find.Execute("91-16=75", $true, $true, $false, $false, $false, $true, 1, $false, "32-11=21", 2) | Out-Null
$d.Content.Find.Execute("77-62=15", $true, $true, $false, $false, $false, $true, 1, $false, "53-2=51", 2) | Out-Null
$d.Content.Find.Execute("92-8=84", $true, $true, $false, $false, $false, $true, 1, $false, "57+39=96", 2) | Out-Null
$d.Content.Find.Execute("51-13=38", $true, $true, $false, $false, $false, $true, 1, $false, "55+2=57", 2) | Out-Null
$d.Content.Find.Execute("59+2=61", $true, $true, $false, $false, $false, $true, 1, $false, "68-23=45", 2) | Out-Null
$d.Content.Find.Execute("39-27=12", $true, $true, $false, $false, $false, $true, 1, $false, "18+70=88", 2) | Out-Null
$d.Content.Find.Execute("34+55=89", $true, $true, $false, $false, $false, $true, 1, $false, "79-33=46", 2) | Out-Null
$d.Content.Find.Execute("32-2=30", $true, $true, $false, $false, $false, $true, 1, $false, "75+9=84", 2) | Out-Null
$d.Content.Find.Execute("78-43=35", $true, $true, $false, $false, $false, $true, 1, $false, "71+12=83", 2) | Out-Null
$d.Content.Find.Execute("71-24=47", $true, $true, $false, $false, $false, $true, 1, $false, "51-20=31", 2) | Out-Null
$d.Content.Find.Execute("43-23=20", $true, $true, $false, $false, $false, $true, 1, $false, "10+11=21", 2) | Out-Null
$d.Content.Find.Execute("93-17=76", $true, $true, $false, $false, $false, $true, 1, $false, "82+17=99", 2) | Out-Null
$d.Content.Find.Execute("57-44=13", $true, $true, $false, $false, $false, $true, 1, $false, "43+6=49", 2) | Out-Null
$d.Content.Find.Execute("81-17=64", $true, $true, $false, $false, $false, $true, 1, $false, "70+2=72", 2) | Out-Null
$d.Content.Find.Execute("54-19=35", $true, $true, $false, $false, $false, $true, 1, $false, "47+21=68", 2) | Out-Null
$d.Content.Find.Execute("51-14=37", $true, $true, $false, $false, $false, $true, 1, $false, "69-60=9", 2) | Out-Null
$d.Content.Find.Execute("91-46=45", $true, $true, $false, $false, $false, $true, 1, $false, "81-79=2", 2) | Out-Null
$d.Content.Find.Execute("55+11=66", $true, $true, $false, $false, $false, $true, 1, $false, "93-72=21", 2) | Out-Null
$d.Content.Find.Execute("85-36=49", $true, $true, $false, $false, $false, $true, 1, $false, "35+15=50", 2) | Out-Null
$d.Content.Find.Execute("64-9=55", $true, $true, $false, $false, $false, $true, 1, $false, "47+52=99", 2) | Out-Null
$d.Content.Find.Execute("56-56=0", $true, $true, $false, $false, $false, $true, 1, $false, "8+34=42", 2) | Out-Null
$d.Content.Find.Execute("11+42=53", $true, $true, $false, $false, $false, $true, 1, $false, "46+2=48", 2) | Out-Null
$d.Content.Find.Execute("52+31=83", $true, $true, $false, $false, $false, $true, 1, $false, "21-2=19", 2) | Out-Null
$d.Content.Find.Execute("51+46=97", $true, $true, $false, $false, $false, $true, 1, $false, "48+49=97", 2) | Out-Null
$d.Content.Find.Execute("26+54=80", $true, $true, $false, $false, $false, $true, 1, $false, "3+91=94", 2) | Out-Null
$d.Content.Find.Execute("23+69=92", $true, $true, $false, $false, $false, $true, 1, $false, "23-20=3", 2) | Out-Null
$d.Content.Find.Execute("89-84=5", $true, $true, $false, $false, $false, $true, 1, $false, "15+67=82", 2) | Out-Null
$d.Content.Find.Execute("58+14=72", $true, $true, $false, $false, $false, $true, 1, $false, "20+1=21", 2) | Out-Null
$d.Content.Find.Execute("51-27=24", $true, $true, $false, $false, $false, $true, 1, $false, "20+13=33", 2) | Out-Null
$d.Content.Find.Execute("37+21=58", $true, $true, $false, $false, $false, $true, 1, $false, "48-19=29", 2) | Out-Null
$d.Content.Find.Execute("12+16=28", $true, $true, $false, $false, $false, $true, 1, $false, "33+11=44", 2) | Out-Null
$d.Content.Find.Execute("87-22=65", $true, $true, $false, $false, $false, $true, 1, $false, "73+6=79", 2) | Out-Null
$d.Content.Find.Execute("13+36=49", $true, $true, $false, $false, $false, $true, 1, $false, "34+56=90", 2) | Out-Null
$d.Content.Find.Execute("67-38=29", $true, $true, $false, $false, $false, $true, 1, $false, "76-14=62", 2) | Out-Null
